$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "S11" record (Excel row 36). Deleting the entire row shifts
# every row below it up by one and Excel automatically drops the now-unused
# "S11" shared string and renumbers the dimension/sortState ranges.
$ws.Rows.Item(36).Delete()

# The freeze pane that was anchored at column B is no longer part of the
# saved view state after this edit, and the selection moves onto the row
# that now occupies position 36 (previously row 37).
$excel.ActiveWindow.FreezePanes = $false
$ws.Rows.Item(36).Select()
